# Applies the "solved booking duplication problem" edit:
#  - strips stray w:proofErr spell-check markers and merges the runs
#    they used to separate (the text/formatting itself is unchanged)
#  - inserts one new bullet paragraph:
#      "- dodać sumowanie ceny i wyświetlenie w podsumowaniu."
#    right before the "- dopracować template ..." paragraph.

$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParagraphXml($findText, $xmlBody) {
    $rng = $d.Content
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Paragraph anchor text not found: $findText"
    }
    $para = $rng.Paragraphs(1).Range
    $para.InsertXML($xmlBody)
}

# --- paragraph 1: "Topic ..." -----------------------------------------
Set-ParagraphXml "Topic" (
    '<w:p ' + $wns + '><w:r><w:t>Topic – główny temat to będzie Car (samochód, który będzie można wynająć)</w:t></w:r></w:p>'
)

# --- paragraph 2: "Entry ..." -------------------------------------------
Set-ParagraphXml "Entry –" (
    '<w:p ' + $wns + '><w:r><w:t>Entry – to będzie wynajęcie Rent (foreignKey, będzie zawierał usera, datę rozpoczęcia i zakończenia)</w:t></w:r></w:p>'
)

# --- "Forms.Form dla rezerwacji ..." -------------------------------------
Set-ParagraphXml "Forms.Form" (
    '<w:p ' + $wns + '><w:r><w:t>Forms.Form dla rezerwacji i filtrowania – ale przed zapisem porównać z cleaned data z zajęć. Czy wyświetli bez zapisywania na głównej stronie, bez przechodzenia na inny template.</w:t></w:r></w:p>'
)

# --- "- zapisać rezerwację, bez podania ForeignKey-a" (highlighted) ------
Set-ParagraphXml "zapisać rezerwację" (
    '<w:p ' + $wns + '><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>- zapisać rezerwację, bez podania ForeignKey-a</w:t></w:r></w:p>'
)

# --- "- Po naciśnięciu buttona potwierdź, zapisać ForeignKey-a (...)" ----
Set-ParagraphXml "Po naciśnięciu" (
    '<w:p ' + $wns + '><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>- Po naciśnięciu buttona potwierdź, zapisać ForeignKey-a (dla auta, a w przyszłości też dla użytkownika)</w:t></w:r></w:p>'
)

# --- "- car_brand" (paragraph mark also carries the highlight via pPr) --
Set-ParagraphXml "car_brand" (
    '<w:p ' + $wns + '><w:pPr><w:rPr><w:highlight w:val="green"/></w:rPr></w:pPr><w:r><w:tab/></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>- car_brand</w:t></w:r></w:p>'
)

# --- "- car_model" --------------------------------------------------------
Set-ParagraphXml "car_model" (
    '<w:p ' + $wns + '><w:pPr><w:rPr><w:highlight w:val="green"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:tab/><w:t>- car_model</w:t></w:r></w:p>'
)

# --- "- year" ---------------------------------------------------------------
Set-ParagraphXml "year" (
    '<w:p ' + $wns + '><w:pPr><w:rPr><w:highlight w:val="green"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:tab/><w:t>- year</w:t></w:r></w:p>'
)

# --- "- registration_number" -------------------------------------------
Set-ParagraphXml "registration_number" (
    '<w:p ' + $wns + '><w:pPr><w:rPr><w:highlight w:val="green"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:tab/><w:t>- registration_number</w:t></w:r></w:p>'
)

# --- "- seats" --------------------------------------------------------------
Set-ParagraphXml "seats" (
    '<w:p ' + $wns + '><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:tab/><w:t>- seats</w:t></w:r></w:p>'
)

# --- "- transmission (manual or automatic)" -------------------------------
Set-ParagraphXml "transmission" (
    '<w:p ' + $wns + '><w:r><w:tab/></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>- transmission (manual or automatic)</w:t></w:r></w:p>'
)

# --- "- category (small, medium, large, kombi, minivan, SUV" -------------
Set-ParagraphXml "category" (
    '<w:p ' + $wns + '><w:r><w:tab/></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>- category (small, medium, large, kombi, minivan, SUV</w:t></w:r></w:p>'
)

# --- "- fuel (benzyna, LPG, elektryczne)" ---------------------------------
Set-ParagraphXml "fuel" (
    '<w:p ' + $wns + '><w:r><w:tab/></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>- fuel (benzyna, LPG, elektryczne)</w:t></w:r></w:p>'
)

# --- "- price per day" ------------------------------------------------------
Set-ParagraphXml "price" (
    '<w:p ' + $wns + '><w:r><w:tab/></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>- price per day</w:t></w:r></w:p>'
)

# --- insert new paragraph before "- dopracować template ..." and clean it up
Set-ParagraphXml "dopracować template" (
    '<w:p ' + $wns + '><w:r><w:t>- dodać sumowanie ceny i wyświetlenie w podsumowaniu.</w:t></w:r></w:p>' +
    '<w:p ' + $wns + '><w:r><w:t>- dopracować template aby wyświetlało wszystkie potrzebne informacje</w:t></w:r></w:p>'
)

# --- "- dodać style lub boostrap-a" -----------------------------------------
Set-ParagraphXml "boostrap" (
    '<w:p ' + $wns + '><w:r><w:t>- dodać style lub boostrap-a</w:t></w:r></w:p>'
)
